$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.319.05'
$ws.Range("E2").Value = '  +1.71%  '
$ws.Range("D3").Value = '3.073.76'
$ws.Range("E3").Value = '  +0.67%  '
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").Value = '''559.24'
$ws.Range("E5").Value = '  +1.73%  '
$ws.Range("D6").Value = '''145.61'
$ws.Range("E6").Value = '  +5.02%  '
$ws.Range("D7").Value = '''0.998'
$ws.Range("E7").Value = '  -0.19%  '
$ws.Range("D8").Value = '3.071.38'
$ws.Range("E8").Value = '  +0.90%  '
$ws.Range("D9").Value = '''0.502'
$ws.Range("E9").Value = '  +0.63%  '
$ws.Range("E10").Value = '  +2.65%  '
$ws.Range("D11").Value = '''6.13'
$ws.Range("E11").Value = '  -2.15%  '
$ws.Range("D12").Value = '''0.469'
$ws.Range("E12").Value = '  +4.27%  '
$ws.Range("E13").Value = '  +1.17%  '
$ws.Range("D14").Value = '''35.19'
$ws.Range("E14").Value = '  +1.77%  '
$ws.Range("D15").Value = '3.569.71'
$ws.Range("E15").Value = '  +0.51%  '
$ws.Range("D16").Value = '64.261.66'
$ws.Range("E16").Value = '  +1.50%  '
$ws.Range("D17").Value = '3.076.06'
$ws.Range("E17").Value = '  +0.67%  '
$ws.Range("E18").Value = '  +1.32%  '
$ws.Range("D19").Value = '''6.78'
$ws.Range("E19").Value = '  +1.34%  '
$ws.Range("D20").Value = '''479.63'
$ws.Range("E20").Value = '  +0.01%  '
$ws.Range("D21").Value = '''13.91'
$ws.Range("E22").Value = '  +0.72%  '
$ws.Range("D23").Value = '''7.58'
$ws.Range("E23").Value = '  +6.09%  '
$ws.Range("D24").Value = '''13.60'
$ws.Range("E24").Value = '  +9.97%  '
$ws.Range("D25").Value = '''81.82'
$ws.Range("E25").Value = '  +1.14%  '
$ws.Range("D26").Value = '''0.999'
$ws.Range("E26").Value = '  -0.07%  '
$ws.Range("E27").Value = '  +2.30%  '
$ws.Range("D28").Value = '''8.12'
$ws.Range("E28").Value = '  +3.50%  '
$ws.Range("E29").Value = '  +5.15%  '
$ws.Range("D30").Value = '''0.999'
$ws.Range("E30").Value = '  -0.06%  '
$ws.Range("D31").Value = '''26.23'
$ws.Range("E32").Value = '  +0.46%  '
$ws.Range("E33").Value = '  +3.77%  '
$ws.Range("D34").Value = '''5.59'
$ws.Range("E34").Value = '  -0.95%  '
$ws.Range("D35").Value = '''6.19'
$ws.Range("E35").Value = '  +3.97%  '
$ws.Range("D36").Value = '''54.93'
$ws.Range("E36").Value = '  -1.09%  '
$ws.Range("D37").Value = '''460.23'
$ws.Range("E37").Value = '  -0.31%  '
$ws.Range("D38").Value = '''3.01'
$ws.Range("E38").Value = '  +18.34%  '
$ws.Range("D39").Value = '''0.0830'
$ws.Range("E39").Value = '  +2.44%  '
$ws.Range("E40").Value = '  +3.56%  '
$ws.Range("D41").Value = '2.974.11'
$ws.Range("E41").Value = '  -4.55%  '
$ws.Range("D42").Value = '''8.27'
$ws.Range("E42").Value = '  +0.82%  '
$ws.Range("E43").Value = '  -2.10%  '
$ws.Range("E44").Value = '  -0.24%  '
$ws.Range("E45").Value = '  +4.89%  '
$ws.Range("D46").Value = '''2.15'
$ws.Range("E46").Value = '  +5.92%  '
$ws.Range("E47").Value = '  +0.07%  '
$ws.Range("E48").Value = '  +2.76%  '
$ws.Range("D49").Value = '''120.72'
$ws.Range("E49").Value = '  +4.24%  '
$ws.Range("E50").Value = '  +2.52%  '
$ws.Range("D51").Value = '''2.08'
$ws.Range("E51").Value = '  +1.37%  '
